$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ----- Row 50: fill in the new "Exp 51" experiment-parameter columns A-D -----
$ws.Range("A50").Value = "Exp 51"
$ws.Range("B50").Value = 0.2
$ws.Range("C50").Value = 1
$ws.Range("D50").Value = "Local"

# Match the centered style already used by the neighbouring A:E cells (row 49)
$ws.Range("A50:D50").Style = $ws.Range("A49:D49").Style
$ws.Range("A50:D50").HorizontalAlignment = $ws.Range("A49:D49").HorizontalAlignment

# New columns O (Test Image) / P (O/P) on row 50
$ws.Range("O50").Value = "Test Image"
$ws.Range("P50").Value = "O/P"

# ----- Row 51: brand-new results row for "Exp 51" -----
$ws.Range("F51").Value = "9,1,0"
$ws.Range("F51").Style = $ws.Range("F50").Style
$ws.Range("F51").HorizontalAlignment = $ws.Range("F50").HorizontalAlignment

$ws.Range("G51").Value = "Exp 51"
$ws.Range("H51").Value = 28
$ws.Range("I51").Value = 64
$ws.Range("J51").Value = 76.63
$ws.Range("K51").Value = 70.51
$ws.Range("L51").Value = 72.69
$ws.Range("M51").Value = 56.04
$ws.Range("N51").Value = 54.29
$ws.Range("O51").Value = 5
$ws.Range("P51").Value = "N/A"
$ws.Range("Q51").Value = "Good - No Overlapping between micro and macro "

# ----- Column widths for the newly used O and Q columns -----
$ws.Columns.Item(15).ColumnWidth = 12.109375
$ws.Columns.Item(17).ColumnWidth = 45.44140625

# ----- View: scroll/select like the author left it -----
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("Q1:Q1048576").Select()
$ws.Range("Q40").Activate()
